$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 376.66
$ws.Range("F2").Value = 1.26
$ws.Range("G2").Value = 0.1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1.13
$ws.Range("O2").Value = 0.5

# Row 3
$ws.Range("E3").Value = 1.04
$ws.Range("F3").Value = 159.74
$ws.Range("G3").Value = 0.1
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 1.26
$ws.Range("O3").Value = 0.5

# Row 4
$ws.Range("E4").Value = 1.14
$ws.Range("F4").Value = 6.74
$ws.Range("G4").Value = 0.1
$ws.Range("M4").Value = 0.59
$ws.Range("N4").Value = 0.36
$ws.Range("O4").Value = 0.5

# Row 5
$ws.Range("E5").Value = 1.39
$ws.Range("F5").Value = 1.46
$ws.Range("G5").Value = 0.1
$ws.Range("M5").Value = 2354.37
$ws.Range("N5").Value = 10.93
$ws.Range("O5").Value = 113.65

# Row 6
$ws.Range("E6").Value = 1.37
$ws.Range("F6").Value = 1.1
$ws.Range("G6").Value = 0.1
$ws.Range("M6").Value = 35263.31
$ws.Range("N6").Value = 332.04
$ws.Range("O6").Value = 14197.41
$ws.Range("T6").Value = 0.1
